$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 39
$ws.Cells.Item(39, 8).Value = 153.3
$ws.Cells.Item(39, 9).Value = 52.5
$ws.Cells.Item(39, 10).Value = 304.5
$ws.Cells.Item(39, 11).Value = 157.5
$ws.Cells.Item(39, 12).Value = 913.5
$ws.Cells.Item(39, 13).Value = 138.5
$ws.Cells.Item(39, 14).Value = -1505.5
# Row 95
$ws.Cells.Item(95, 8).Value = 39800
$ws.Cells.Item(95, 10).Value = 39800
$ws.Cells.Item(95, 12).Value = 39800
$ws.Cells.Item(95, 14).Value = -45292
# Row 105
$ws.Cells.Item(105, 8).Value = 45221
$ws.Cells.Item(105, 10).Value = 45221
$ws.Cells.Item(105, 12).Value = 45221
$ws.Cells.Item(105, 14).Value = -52209
# Row 108
$ws.Cells.Item(108, 8).Value = 34074.5
$ws.Cells.Item(108, 10).Value = 34074.5
$ws.Cells.Item(108, 12).Value = 34074.5
$ws.Cells.Item(108, 14).Value = -41754.5
# Row 109
$ws.Cells.Item(109, 8).Value = 34219
$ws.Cells.Item(109, 10).Value = 34219
$ws.Cells.Item(109, 12).Value = 34219
$ws.Cells.Item(109, 14).Value = -36993
# Row 114
$ws.Cells.Item(114, 8).Value = 45718
$ws.Cells.Item(114, 10).Value = 45718
$ws.Cells.Item(114, 12).Value = 45718
$ws.Cells.Item(114, 14).Value = -54396
# Row 120
$ws.Cells.Item(120, 8).Value = 49706
$ws.Cells.Item(120, 10).Value = 49706
$ws.Cells.Item(120, 12).Value = 49706
$ws.Cells.Item(120, 14).Value = -59382
# Row 123
$ws.Cells.Item(123, 8).Value = 32222.133
$ws.Cells.Item(123, 10).Value = 32222.133
$ws.Cells.Item(123, 12).Value = 32222.133
$ws.Cells.Item(123, 14).Value = -42022.133
# Row 124
$ws.Cells.Item(124, 8).Value = 48078.75
$ws.Cells.Item(124, 10).Value = 48078.75
$ws.Cells.Item(124, 12).Value = 48078.75
$ws.Cells.Item(124, 14).Value = -57898.75
# Row 126
$ws.Cells.Item(126, 8).Value = 46764
$ws.Cells.Item(126, 10).Value = 46764
$ws.Cells.Item(126, 12).Value = 46764
$ws.Cells.Item(126, 14).Value = -56644
# Row 130
$ws.Cells.Item(130, 8).Value = 46501.332
$ws.Cells.Item(130, 10).Value = 46501.332
$ws.Cells.Item(130, 12).Value = 46501.332
$ws.Cells.Item(130, 14).Value = -56541.332

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 2072.5833
$ws.Cells.Item(2, 9).Value = 2155
$ws.Cells.Item(2, 10).Value = 1825.3334
$ws.Cells.Item(2, 11).Value = 2155
$ws.Cells.Item(2, 12).Value = 1825.3334
$ws.Cells.Item(2, 13).Value = -2042
$ws.Cells.Item(2, 14).Value = -2051.3334
# Row 59
$ws.Cells.Item(59, 8).Value = 41998
$ws.Cells.Item(59, 10).Value = 41998
$ws.Cells.Item(59, 12).Value = 41998
$ws.Cells.Item(59, 14).Value = -43606
# Row 74
$ws.Cells.Item(74, 8).Value = 1979.4286
$ws.Cells.Item(74, 9).Value = 905.4545000000001
$ws.Cells.Item(74, 10).Value = 3160.8
$ws.Cells.Item(74, 11).Value = 905.4545000000001
$ws.Cells.Item(74, 12).Value = 3160.8
$ws.Cells.Item(74, 13).Value = -31.45450000000005
$ws.Cells.Item(74, 14).Value = -4908.8
# Row 77
$ws.Cells.Item(77, 8).Value = 1979.4286
$ws.Cells.Item(77, 9).Value = 905.4545000000001
$ws.Cells.Item(77, 10).Value = 3160.8
$ws.Cells.Item(77, 11).Value = 4527.2725
$ws.Cells.Item(77, 12).Value = 15804
$ws.Cells.Item(77, 13).Value = -159.2725
$ws.Cells.Item(77, 14).Value = -24540
# Row 116
$ws.Cells.Item(116, 8).Value = 2072.5833
$ws.Cells.Item(116, 9).Value = 2155
$ws.Cells.Item(116, 10).Value = 1825.3334
$ws.Cells.Item(116, 11).Value = 2155
$ws.Cells.Item(116, 12).Value = 1825.3334
$ws.Cells.Item(116, 13).Value = 139
$ws.Cells.Item(116, 14).Value = -6413.3334
# Row 135
$ws.Cells.Item(135, 8).Value = 36377.5
$ws.Cells.Item(135, 10).Value = 36377.5
$ws.Cells.Item(135, 12).Value = 36377.5
$ws.Cells.Item(135, 14).Value = -46517.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 2072.5833
$ws.Cells.Item(3, 9).Value = 2155
$ws.Cells.Item(3, 10).Value = 1825.3334
$ws.Cells.Item(3, 11).Value = 2155
$ws.Cells.Item(3, 12).Value = 1825.3334
$ws.Cells.Item(3, 13).Value = -2041
$ws.Cells.Item(3, 14).Value = -2053.3334
# Row 5
$ws.Cells.Item(5, 8).Value = 6866
$ws.Cells.Item(5, 9).Value = 299
$ws.Cells.Item(5, 11).Value = 299
$ws.Cells.Item(5, 13).Value = -186
# Row 44
$ws.Cells.Item(44, 8).Value = 50000
$ws.Cells.Item(44, 10).Value = 50000
$ws.Cells.Item(44, 12).Value = 50000
$ws.Cells.Item(44, 14).Value = -50994
# Row 132
$ws.Cells.Item(132, 8).Value = 39720
$ws.Cells.Item(132, 10).Value = 39720
$ws.Cells.Item(132, 12).Value = 39720
$ws.Cells.Item(132, 14).Value = -49840

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 75
$ws.Cells.Item(75, 8).Value = 19285.715
# Row 78
$ws.Cells.Item(78, 8).Value = 19285.715
# Row 92
$ws.Cells.Item(92, 8).Value = 35187.547
$ws.Cells.Item(92, 10).Value = 35187.547
$ws.Cells.Item(92, 12).Value = 35187.547
$ws.Cells.Item(92, 14).Value = -40179.547

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Cells.Item(22, 8).Value = 950
$ws.Cells.Item(22, 9).Value = 600
$ws.Cells.Item(22, 10).Value = 1300
$ws.Cells.Item(22, 11).Value = 1800
$ws.Cells.Item(22, 12).Value = 3900
$ws.Cells.Item(22, 13).Value = -1631
$ws.Cells.Item(22, 14).Value = -4238
# Row 27
$ws.Cells.Item(27, 8).Value = 950
$ws.Cells.Item(27, 9).Value = 600
$ws.Cells.Item(27, 10).Value = 1300
$ws.Cells.Item(27, 11).Value = 1800
$ws.Cells.Item(27, 12).Value = 3900
$ws.Cells.Item(27, 13).Value = -1698
$ws.Cells.Item(27, 14).Value = -4104
# Row 131
$ws.Cells.Item(131, 8).Value = 3733.439
$ws.Cells.Item(131, 9).Value = 7105.067
$ws.Cells.Item(131, 10).Value = 1788.2693
$ws.Cells.Item(131, 11).Value = 21315.201
$ws.Cells.Item(131, 12).Value = 5364.8079
$ws.Cells.Item(131, 13).Value = -16275.201
$ws.Cells.Item(131, 14).Value = -15444.8079
# Row 140
$ws.Cells.Item(140, 8).Value = 2764.1428
$ws.Cells.Item(140, 9).Value = 1888.8889
$ws.Cells.Item(140, 10).Value = 4339.6
$ws.Cells.Item(140, 11).Value = 5666.6667
$ws.Cells.Item(140, 12).Value = 13018.8
$ws.Cells.Item(140, 13).Value = -486.6666999999998
$ws.Cells.Item(140, 14).Value = -23378.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 47
$ws.Cells.Item(47, 8).Value = 8000
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 11).Value = 0
$ws.Cells.Item(47, 13).ClearContents()
# Row 105
$ws.Cells.Item(105, 8).Value = 42219
$ws.Cells.Item(105, 10).Value = 42219
$ws.Cells.Item(105, 12).Value = 42219
$ws.Cells.Item(105, 14).Value = -49207
# Row 110
$ws.Cells.Item(110, 8).Value = 47698
$ws.Cells.Item(110, 10).Value = 47698
$ws.Cells.Item(110, 12).Value = 47698
$ws.Cells.Item(110, 14).Value = -55878
# Row 120
$ws.Cells.Item(120, 8).Value = 38301
$ws.Cells.Item(120, 10).Value = 38301
$ws.Cells.Item(120, 12).Value = 38301
$ws.Cells.Item(120, 14).Value = -47977
# Row 130
$ws.Cells.Item(130, 8).Value = 46135.855
$ws.Cells.Item(130, 10).Value = 46135.855
$ws.Cells.Item(130, 12).Value = 46135.855
$ws.Cells.Item(130, 14).Value = -56175.855
# Row 132
$ws.Cells.Item(132, 8).Value = 3693.2285
$ws.Cells.Item(132, 9).Value = 1331.8572
$ws.Cells.Item(132, 11).Value = 3995.5716
$ws.Cells.Item(132, 13).Value = -1465.5716

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 48
$ws.Cells.Item(48, 8).Value = 26360
$ws.Cells.Item(48, 10).Value = 26360
$ws.Cells.Item(48, 12).Value = 26360
$ws.Cells.Item(48, 14).Value = -27498
# Row 127
$ws.Cells.Item(127, 8).Value = 15379.111
$ws.Cells.Item(127, 10).Value = 42412
$ws.Cells.Item(127, 12).Value = 42412
$ws.Cells.Item(127, 14).Value = -52332
